$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New chronological (ascending) order of "Periodo Mora" labels for rows 16-39.
# Previously the periods were listed in descending order (2111 down to 1912);
# the database was corrected to ascending order (1912 up to 2111).
$periods = @("1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# The "Valor Mora" (F column) value that used to sit with the oldest period
# row (row 16) now belongs with the newest period (row 39), and vice versa.
$f16 = $ws.Cells.Item(16, 6).Value2
$f39 = $ws.Cells.Item(39, 6).Value2
$ws.Cells.Item(16, 6).Value = $f39
$ws.Cells.Item(39, 6).Value = $f16
